$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 259; everything from old row 259 downward shifts to 260..275
$ws.Rows.Item(259).Insert()

# Populate the newly inserted row 259 with the new weekly record
$ws.Cells.Item(259, 1).Value  = 8
$ws.Cells.Item(259, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(259, 3).Value  = "Coquimbo"
$ws.Cells.Item(259, 4).Value  = 44516
$ws.Cells.Item(259, 5).Value  = 4
$ws.Cells.Item(259, 6).Value  = 100114001
$ws.Cells.Item(259, 7).Value  = "Papa"
$ws.Cells.Item(259, 8).Value  = "Cardinal"
$ws.Cells.Item(259, 9).Value  = "1a nueva(o)"
$ws.Cells.Item(259, 10).Value = 2600
$ws.Cells.Item(259, 11).Value = 11500
$ws.Cells.Item(259, 12).Value = 12000
$ws.Cells.Item(259, 13).Value = 11750
$ws.Cells.Item(259, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(259, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(259, 16).Value = 470
$ws.Cells.Item(259, 17).Value = 25
$ws.Cells.Item(259, 18).Value = "Hortaliza"
